# Auto-generated Excel COM-interop script applying profit/price refresh
# from the scheduled market-data runner (Sargatanas_Profits workbook).
$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("ALC")
$ws.Cells.Item(87, 8).Value = 54519.4
$ws.Cells.Item(87, 10).Value = 54519.4
$ws.Cells.Item(87, 12).Value = 54519.4
$ws.Cells.Item(87, 14).Value = -57015.4

$ws.Cells.Item(90, 8).Value = 54519.4
$ws.Cells.Item(90, 10).Value = 54519.4
$ws.Cells.Item(90, 12).Value = 163558.2
$ws.Cells.Item(90, 14).Value = -176038.2

$ws.Cells.Item(112, 8).Value = 4245.8076
$ws.Cells.Item(112, 10).Value = 4832
$ws.Cells.Item(112, 12).Value = 14496
$ws.Cells.Item(112, 14).Value = -16712

$ws.Cells.Item(132, 8).Value = 809.0192
$ws.Cells.Item(132, 9).Value = 809.0192
$ws.Cells.Item(132, 11).Value = 2427.0576
$ws.Cells.Item(132, 13).Value = 102.9423999999999

$ws = $wb.Worksheets.Item("ARM")
$ws.Cells.Item(32, 8).Value = 3231772.8
$ws.Cells.Item(32, 9).Value = 3283933
$ws.Cells.Item(32, 11).Value = 3283933
$ws.Cells.Item(32, 13).Value = -3283646

$ws.Cells.Item(74, 8).Value = 2745.8572
$ws.Cells.Item(74, 9).Value = 1781.1333
$ws.Cells.Item(74, 11).Value = 1781.1333
$ws.Cells.Item(74, 13).Value = -907.1333

$ws.Cells.Item(77, 8).Value = 2745.8572
$ws.Cells.Item(77, 9).Value = 1781.1333
$ws.Cells.Item(77, 11).Value = 8905.666499999999
$ws.Cells.Item(77, 13).Value = -4537.666499999999

$ws.Cells.Item(97, 8).Value = 5954571
$ws.Cells.Item(97, 9).Value = 2215.9
$ws.Cells.Item(97, 11).Value = 2215.9
$ws.Cells.Item(97, 13).Value = -1719.9

$ws.Cells.Item(110, 8).Value = 55556976
$ws.Cells.Item(110, 9).Value = 1700
$ws.Cells.Item(110, 10).Value = 111112250
$ws.Cells.Item(110, 11).Value = 1700
$ws.Cells.Item(110, 12).Value = 111112250
$ws.Cells.Item(110, 13).Value = 345
$ws.Cells.Item(110, 14).Value = -111116340

$ws.Cells.Item(122, 8).Value = 3661.9167
$ws.Cells.Item(122, 9).Value = 2937.375
$ws.Cells.Item(122, 11).Value = 8812.125
$ws.Cells.Item(122, 13).Value = -6362.125

$ws.Cells.Item(130, 8).Value = 26975.666
$ws.Cells.Item(130, 10).Value = 26975.666
$ws.Cells.Item(130, 12).Value = 26975.666
$ws.Cells.Item(130, 14).Value = -37015.666

$ws.Cells.Item(132, 8).Value = 5624.5405
$ws.Cells.Item(132, 9).Value = 1781.45
$ws.Cells.Item(132, 10).Value = 10145.823
$ws.Cells.Item(132, 11).Value = 5344.35
$ws.Cells.Item(132, 12).Value = 30437.469
$ws.Cells.Item(132, 13).Value = -2814.35
$ws.Cells.Item(132, 14).Value = -35497.469

$ws.Cells.Item(133, 8).Value = 122222
$ws.Cells.Item(133, 10).Value = 122222
$ws.Cells.Item(133, 12).Value = 122222
$ws.Cells.Item(133, 14).Value = -127282

$ws = $wb.Worksheets.Item("BSM")
$ws.Cells.Item(80, 8).Value = 21739534
$ws.Cells.Item(80, 9).Value = 50000388
$ws.Cells.Item(80, 11).Value = 50000388
$ws.Cells.Item(80, 13).Value = -49999390

$ws.Cells.Item(83, 8).Value = 21739534
$ws.Cells.Item(83, 9).Value = 50000388
$ws.Cells.Item(83, 11).Value = 250001940
$ws.Cells.Item(83, 13).Value = -249996948

$ws.Cells.Item(107, 8).Value = 375000500
$ws.Cells.Item(107, 10).Value = 1500
$ws.Cells.Item(107, 12).Value = 1500
$ws.Cells.Item(107, 14).Value = -5340

$ws = $wb.Worksheets.Item("CRP")
$ws.Cells.Item(31, 8).Value = 7328.4443
$ws.Cells.Item(31, 9).Value = 2527.3333
$ws.Cells.Item(31, 10).Value = 10529.186
$ws.Cells.Item(31, 11).Value = 2527.3333
$ws.Cells.Item(31, 12).Value = 10529.186
$ws.Cells.Item(31, 13).Value = -2232.3333
$ws.Cells.Item(31, 14).Value = -11119.186

$ws.Cells.Item(34, 8).Value = 7328.4443
$ws.Cells.Item(34, 9).Value = 2527.3333
$ws.Cells.Item(34, 10).Value = 10529.186
$ws.Cells.Item(34, 11).Value = 2527.3333
$ws.Cells.Item(34, 12).Value = 10529.186
$ws.Cells.Item(34, 13).Value = -2325.3333
$ws.Cells.Item(34, 14).Value = -10933.186

$ws.Cells.Item(76, 8).Value = 4718.5
$ws.Cells.Item(76, 9).Value = 4718.5
$ws.Cells.Item(76, 11).Value = 4718.5
$ws.Cells.Item(76, 13).Value = -4403.5

$ws.Cells.Item(79, 8).Value = 4718.5
$ws.Cells.Item(79, 9).Value = 4718.5
$ws.Cells.Item(79, 11).Value = 4718.5
$ws.Cells.Item(79, 13).Value = -3626.5

$ws.Cells.Item(86, 8).Value = 14886275
$ws.Cells.Item(86, 10).Value = 4147.25
$ws.Cells.Item(86, 12).Value = 4147.25
$ws.Cells.Item(86, 14).Value = -6393.25

$ws.Cells.Item(89, 8).Value = 14886275
$ws.Cells.Item(89, 10).Value = 4147.25
$ws.Cells.Item(89, 12).Value = 20736.25
$ws.Cells.Item(89, 14).Value = -31968.25

$ws.Cells.Item(132, 8).Value = 6262.3657
$ws.Cells.Item(132, 9).Value = 3798.8635
$ws.Cells.Item(132, 10).Value = 9114.842000000001
$ws.Cells.Item(132, 11).Value = 11396.5905
$ws.Cells.Item(132, 12).Value = 27344.526
$ws.Cells.Item(132, 13).Value = -8866.5905
$ws.Cells.Item(132, 14).Value = -32404.526

$ws.Cells.Item(134, 8).Value = 6466.778
$ws.Cells.Item(134, 9).Value = 2251
$ws.Cells.Item(134, 11).Value = 6753
$ws.Cells.Item(134, 13).Value = -4218

$ws = $wb.Worksheets.Item("CUL")
$ws.Cells.Item(4, 8).Value = 50855580
$ws.Cells.Item(4, 9).Value = 57849868
$ws.Cells.Item(4, 11).Value = 173549604
$ws.Cells.Item(4, 13).Value = -173549492

$ws.Cells.Item(32, 8).Value = 108.166664
$ws.Cells.Item(32, 9).Value = 99
$ws.Cells.Item(32, 11).Value = 297
$ws.Cells.Item(32, 13).Value = -14

$ws.Cells.Item(80, 8).Value = 31254476
$ws.Cells.Item(80, 9).Value = 23813392
$ws.Cells.Item(80, 10).Value = 45460180
$ws.Cells.Item(80, 11).Value = 71440176
$ws.Cells.Item(80, 12).Value = 136380540
$ws.Cells.Item(80, 13).Value = -71439240
$ws.Cells.Item(80, 14).Value = -136382412

$ws.Cells.Item(83, 8).Value = 31254476
$ws.Cells.Item(83, 9).Value = 23813392
$ws.Cells.Item(83, 10).Value = 45460180
$ws.Cells.Item(83, 11).Value = 214320528
$ws.Cells.Item(83, 12).Value = 409141620
$ws.Cells.Item(83, 13).Value = -214315848
$ws.Cells.Item(83, 14).Value = -409150980

$ws.Cells.Item(98, 8).Value = 3379.8
$ws.Cells.Item(98, 10).Value = 3165.6667
$ws.Cells.Item(98, 12).Value = 9497.000100000001
$ws.Cells.Item(98, 14).Value = -12493.0001

$ws.Cells.Item(113, 8).Value = 5682.278
$ws.Cells.Item(113, 9).Value = 1438.625
$ws.Cells.Item(113, 10).Value = 9077.200000000001
$ws.Cells.Item(113, 11).Value = 4315.875
$ws.Cells.Item(113, 12).Value = 27231.6
$ws.Cells.Item(113, 13).Value = -2145.875
$ws.Cells.Item(113, 14).Value = -31571.6

$ws.Cells.Item(133, 8).Value = 11007.115
$ws.Cells.Item(133, 9).Value = 5013.4287
$ws.Cells.Item(133, 11).Value = 15040.2861
$ws.Cells.Item(133, 13).Value = -9980.286100000001

$ws = $wb.Worksheets.Item("GSM")
$ws.Cells.Item(2, 8).Value = 1176683.2
$ws.Cells.Item(2, 9).Value = 102
$ws.Cells.Item(2, 10).Value = 2222533.2
$ws.Cells.Item(2, 11).Value = 102
$ws.Cells.Item(2, 12).Value = 2222533.2
$ws.Cells.Item(2, 13).Value = 11
$ws.Cells.Item(2, 14).Value = -2222759.2

$ws.Cells.Item(80, 8).Value = 2026.5454
$ws.Cells.Item(80, 9).Value = 2366.3333
$ws.Cells.Item(80, 10).Value = 1618.8
$ws.Cells.Item(80, 11).Value = 2366.3333
$ws.Cells.Item(80, 12).Value = 1618.8
$ws.Cells.Item(80, 13).Value = -1368.3333
$ws.Cells.Item(80, 14).Value = -3614.8

$ws.Cells.Item(83, 8).Value = 2026.5454
$ws.Cells.Item(83, 9).Value = 2366.3333
$ws.Cells.Item(83, 10).Value = 1618.8
$ws.Cells.Item(83, 11).Value = 11831.6665
$ws.Cells.Item(83, 12).Value = 8094
$ws.Cells.Item(83, 13).Value = -6839.666499999999
$ws.Cells.Item(83, 14).Value = -18078

$ws.Cells.Item(132, 8).Value = 5202.5713
$ws.Cells.Item(132, 9).Value = 2052.9412
$ws.Cells.Item(132, 10).Value = 10070.182
$ws.Cells.Item(132, 11).Value = 6158.823600000001
$ws.Cells.Item(132, 12).Value = 30210.546
$ws.Cells.Item(132, 13).Value = -3628.823600000001
$ws.Cells.Item(132, 14).Value = -35270.546

$ws = $wb.Worksheets.Item("LTW")
$ws.Cells.Item(16, 8).Value = 723.75
$ws.Cells.Item(16, 9).Value = 782
$ws.Cells.Item(16, 10).Value = 587.8333
$ws.Cells.Item(16, 11).Value = 782
$ws.Cells.Item(16, 12).Value = 587.8333
$ws.Cells.Item(16, 13).Value = -612
$ws.Cells.Item(16, 14).Value = -927.8333

$ws.Cells.Item(22, 8).Value = 1457.6364

$ws.Cells.Item(27, 8).Value = 1457.6364

$ws.Cells.Item(40, 8).Value = 2957.1025
$ws.Cells.Item(40, 9).Value = 2017.6333
$ws.Cells.Item(40, 11).Value = 2017.6333
$ws.Cells.Item(40, 13).Value = -1881.6333

$ws.Cells.Item(68, 8).Value = 6413.5713
$ws.Cells.Item(68, 9).Value = 4278.2
$ws.Cells.Item(68, 10).Value = 7599.8887
$ws.Cells.Item(68, 11).Value = 4278.2
$ws.Cells.Item(68, 12).Value = 7599.8887
$ws.Cells.Item(68, 13).Value = -3529.2
$ws.Cells.Item(68, 14).Value = -9097.8887

$ws.Cells.Item(71, 8).Value = 6413.5713
$ws.Cells.Item(71, 9).Value = 4278.2
$ws.Cells.Item(71, 10).Value = 7599.8887
$ws.Cells.Item(71, 11).Value = 21391
$ws.Cells.Item(71, 12).Value = 37999.4435
$ws.Cells.Item(71, 13).Value = -17647
$ws.Cells.Item(71, 14).Value = -45487.4435

$ws.Cells.Item(82, 8).Value = 521217.28
$ws.Cells.Item(82, 9).Value = 794493.4399999999
$ws.Cells.Item(82, 10).Value = 1992.6
$ws.Cells.Item(82, 11).Value = 794493.4399999999
$ws.Cells.Item(82, 12).Value = 1992.6
$ws.Cells.Item(82, 13).Value = -794132.4399999999
$ws.Cells.Item(82, 14).Value = -2714.6

$ws.Cells.Item(85, 8).Value = 521217.28
$ws.Cells.Item(85, 9).Value = 794493.4399999999
$ws.Cells.Item(85, 10).Value = 1992.6
$ws.Cells.Item(85, 11).Value = 794493.4399999999
$ws.Cells.Item(85, 12).Value = 1992.6
$ws.Cells.Item(85, 13).Value = -793245.4399999999
$ws.Cells.Item(85, 14).Value = -4488.6

$ws.Cells.Item(97, 8).Value = 0
$ws.Cells.Item(97, 10).Value = 0
$ws.Cells.Item(97, 14).ClearContents()

$ws.Cells.Item(115, 8).Value = 52678.332
$ws.Cells.Item(115, 10).Value = 52678.332
$ws.Cells.Item(115, 12).Value = 52678.332
$ws.Cells.Item(115, 14).Value = -55028.332

$ws = $wb.Worksheets.Item("WVR")
$ws.Cells.Item(62, 8).Value = 140355.86
$ws.Cells.Item(62, 9).Value = 163082
$ws.Cells.Item(62, 11).Value = 163082
$ws.Cells.Item(62, 13).Value = -162458

$ws.Cells.Item(65, 8).Value = 140355.86
$ws.Cells.Item(65, 9).Value = 163082
$ws.Cells.Item(65, 11).Value = 815410
$ws.Cells.Item(65, 13).Value = -812290

$ws.Cells.Item(113, 8).Value = 12177.667
$ws.Cells.Item(113, 9).Value = 24241.092
$ws.Cells.Item(113, 10).Value = 1970.1538
$ws.Cells.Item(113, 11).Value = 72723.276
$ws.Cells.Item(113, 12).Value = 5910.4614
$ws.Cells.Item(113, 13).Value = -70553.276
$ws.Cells.Item(113, 14).Value = -10250.4614

$ws.Cells.Item(132, 8).Value = 9921.263000000001
$ws.Cells.Item(132, 9).Value = 14011.333
$ws.Cells.Item(132, 11).Value = 42033.999
$ws.Cells.Item(132, 13).Value = -39503.999

$ws.Cells.Item(138, 8).Value = 73999.39999999999
$ws.Cells.Item(138, 10).Value = 73999.39999999999
$ws.Cells.Item(138, 12).Value = 73999.39999999999
$ws.Cells.Item(138, 14).Value = -84279.39999999999
